# Remove unnecessary information in templates
# (mine-tag-list.xlsx) — strip the pre-filled values that follow the
# bold labels in A4:A7, clear the "net weight / sample / dust" lines
# (B13:B15), clear the export-date value (B8), and move the selection
# to E6 (matching the saved sheet view in the target workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MINE TAGS")

# --- A4: "       MATERIAL : WOLFRAMITE" -> "       MATERIAL :" ------------
# Runs: [1..7]="       " (plain) [8..28]="MATERIAL : WOLFRAMITE" (bold)
$r = $ws.Range("A4")
$bold = $r.Characters(8, 21)
$bold.Text = "MATERIAL :"
$ws.Range("A4").Characters(8, 10).Font.Bold = $true

# --- A5: "       BUYER: HALCYON INC. " -> "       BUYER: " -----------------
# Runs: [1..7]="       " (plain) [8..27]="BUYER: HALCYON INC. " (bold)
$r = $ws.Range("A5")
$bold = $r.Characters(8, 20)
$bold.Text = "BUYER: "
$ws.Range("A5").Characters(8, 7).Font.Bold = $true

# --- A6: "       LOT N°: WKKIG2303" -> "       LOT N°: " --------------------
# Runs: [1..7]="       " (plain) [8..24]="LOT N°: WKKIG2303" (bold)
$r = $ws.Range("A6")
$bold = $r.Characters(8, 17)
$bold.Text = "LOT N°: "
$ws.Range("A6").Characters(8, 8).Font.Bold = $true

# --- A7: "       iTSCi SHIPMENT NUMBER: KZM/RW/ 0000022" -> "       iTSCi SHIPMENT NUMBER:" --
# Runs: [1..5]="     " (plain) [6..45]="  iTSCi SHIPMENT NUMBER: KZM/RW/ 0000022" (bold)
$r = $ws.Range("A7")
$bold = $r.Characters(6, 40)
$bold.Text = "  iTSCi SHIPMENT NUMBER:"
$ws.Range("A7").Characters(6, 24).Font.Bold = $true

# --- B8: "EXPORT DATE : 13/11/2023" -> "EXPORT DATE :" ---------------------
$ws.Range("B8").Value = "EXPORT DATE :"

# --- B13:B15: clear the net-weight / sample / dust lines -------------------
$ws.Range("B13").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("B15").Value = ""

# --- Selection moves to E6 (matches the saved sheet view) ------------------
$ws.Range("E6").Select()
